$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.5410300942337471
$ws.Range("C2").Value = 0.2768127168325076
$ws.Range("D2").Value = 0.03316564316829584
$ws.Range("F2").Value = 0.6696811536956773
$ws.Range("G2").Value = 0.5099755435820512
$ws.Range("H2").Value = 0.6575008361093424
$ws.Range("K2").Value = 0.2658695971831264
$ws.Range("L2").Value = 0.3018766062606488
$ws.Range("N2").Value = 1.46088788364032
$ws.Range("O2").Value = 2.300653869830867
$ws.Range("B3").Value = 0.4973042949314106
$ws.Range("C3").Value = 0.2790564963711972
$ws.Range("D3").Value = 0.03049942358976665
$ws.Range("F3").Value = 0.6693404871908868
$ws.Range("G3").Value = 0.5111656612577136
$ws.Range("H3").Value = 0.6613721809385282
$ws.Range("K3").Value = 0.2322006871786613
$ws.Range("L3").Value = 0.2905448421806653
$ws.Range("N3").Value = 1.475722630817817
$ws.Range("O3").Value = 2.31104569319848
$ws.Range("B4").Value = 0.4705858802843466
$ws.Range("C4").Value = 0.2805091056868099
$ws.Range("D4").Value = 0.0288480876533086
$ws.Range("F4").Value = 0.6695146128997891
$ws.Range("G4").Value = 0.5122153051873326
$ws.Range("H4").Value = 0.6640097691325337
$ws.Range("K4").Value = 0.2114735398789378
$ws.Range("L4").Value = 0.2837571077002963
$ws.Range("N4").Value = 1.485311072507731
$ws.Range("O4").Value = 2.318639534505877
$ws.Range("B5").Value = 0.4597311740519956
$ws.Range("C5").Value = 0.2811199197046541
$ws.Range("D5").Value = 0.02817160032829236
$ws.Range("F5").Value = 0.6696820274203219
$ws.Range("G5").Value = 0.5127232299080049
$ws.Range("H5").Value = 0.665150199880685
$ws.Range("K5").Value = 0.2030139292224078
$ws.Range("L5").Value = 0.2810338972334137
$ws.Range("N5").Value = 1.489339121317988
$ws.Range("O5").Value = 2.322039290798941
$ws.Range("B6").Value = 0.4579307887912876
$ws.Range("C6").Value = 0.281222484972516
$ws.Range("D6").Value = 0.02805905640322948
$ws.Range("F6").Value = 0.6697156544025731
$ws.Range("G6").Value = 0.5128124133179952
$ws.Range("H6").Value = 0.6653435312176299
$ws.Range("K6").Value = 0.2016084397409799
$ws.Range("L6").Value = 0.280584300820351
$ws.Range("N6").Value = 1.490015265198554
$ws.Range("O6").Value = 2.322622255171879
$ws.Range("B7").Value = 0.4704393542715763
$ws.Range("C7").Value = 0.2805172669273039
$ws.Range("D7").Value = 0.02883897865736174
$ws.Range("F7").Value = 0.669516480062434
$ws.Range("G7").Value = 0.5122218305735302
$ws.Range("H7").Value = 0.6640248837066096
$ws.Range("K7").Value = 0.211359503078711
$ws.Range("L7").Value = 0.2837202079218031
$ws.Range("N7").Value = 1.485364907523101
$ws.Range("O7").Value = 2.318684148875164
$ws.Range("B8").Value = 0.5259270023358908
$ws.Range("C8").Value = 0.2775708375148236
$ws.Range("D8").Value = 0.03224931064976744
$ws.Range("F8").Value = 0.6694841811634191
$ws.Range("G8").Value = 0.5103197026326782
$ws.Range("H8").Value = 0.6587816371743997
$ws.Range("K8").Value = 0.2542721317164478
$ws.Range("L8").Value = 0.2979341804341971
$ws.Range("N8").Value = 1.465903325209055
$ws.Range("O8").Value = 2.303985241918056
$ws.Range("B9").Value = 0.6357385174968897
$ws.Range("C9").Value = 0.2723861727396057
$ws.Range("D9").Value = 0.03882260148906624
$ws.Range("F9").Value = 0.6724598087353399
$ws.Range("G9").Value = 0.5091208697625902
$ws.Range("H9").Value = 0.6505641797605648
$ws.Range("K9").Value = 0.3379736791028449
$ws.Range("L9").Value = 0.3271543049267933
$ws.Range("N9").Value = 1.43154517965387
$ws.Range("O9").Value = 2.284782979635423
$ws.Range("B10").Value = 0.7169996276427639
$ws.Range("C10").Value = 0.2689369569658799
$ws.Range("D10").Value = 0.0435811324608224
$ws.Range("F10").Value = 0.6764968938239235
$ws.Range("G10").Value = 0.5097852493180852
$ws.Range("H10").Value = 0.6457816673298709
$ws.Range("K10").Value = 0.3991749111001468
$ws.Range("L10").Value = 0.3494425665349468
$ws.Range("N10").Value = 1.408619749842089
$ws.Range("O10").Value = 2.276537945114114
$ws.Range("B11").Value = 0.7540885934607218
$ws.Range("C11").Value = 0.2674456142169763
$ws.Range("D11").Value = 0.04573030143276924
$ws.Range("F11").Value = 0.6787350842365925
$ws.Range("G11").Value = 0.51042350288688
$ws.Range("H11").Value = 0.643877688979515
$ws.Range("K11").Value = 0.4269491289409473
$ws.Range("L11").Value = 0.3597601811337228
$ws.Range("N11").Value = 1.398693030687461
$ws.Range("O11").Value = 2.274059643386693
$ws.Range("B12").Value = 0.7681502534720153
$ws.Range("C12").Value = 0.2668920327321649
$ws.Range("D12").Value = 0.04654187640673513
$ws.Range("F12").Value = 0.6796403405045766
$ws.Range("G12").Value = 0.5107135360183719
$ws.Range("H12").Value = 0.6431956952012996
$ws.Range("K12").Value = 0.4374564728812231
$ws.Range("L12").Value = 0.3636928100565626
$ws.Range("N12").Value = 1.395006239970831
$ws.Range("O12").Value = 2.273304066671329
$ws.Range("B13").Value = 0.7651210856266175
$ws.Range("C13").Value = 0.2670107604676657
$ws.Range("D13").Value = 0.04636719072570372
$ws.Range("F13").Value = 0.6794428120485847
$ws.Range("G13").Value = 0.5106489219265029
$ws.Range("H13").Value = 0.6433408410478734
$ws.Range("K13").Value = 0.4351939890490542
$ws.Range("L13").Value = 0.3628447124322918
$ws.Range("N13").Value = 1.395797042809605
$ws.Range("O13").Value = 2.273458660103984
$ws.Range("B14").Value = 0.7552451210867162
$ws.Range("C14").Value = 0.2673998472258683
$ws.Range("D14").Value = 0.04579711588198165
$ws.Range("F14").Value = 0.6788084041359568
$ws.Range("G14").Value = 0.5104463951600877
$ws.Range("H14").Value = 0.6438207996116034
$ws.Range("K14").Value = 0.4278137812023033
$ws.Range("L14").Value = 0.3600832089910995
$ws.Range("N14").Value = 1.398388268168919
$ws.Range("O14").Value = 2.273993816395006
$ws.Range("B15").Value = 0.7491979814293472
$ws.Range("C15").Value = 0.2676396266030965
$ws.Range("D15").Value = 0.04544763199290003
$ws.Range("F15").Value = 0.6784273236386227
$ws.Range("G15").Value = 0.5103286376033651
$ws.Range("H15").Value = 0.6441198655849547
$ws.Range("K15").Value = 0.4232918517236328
$ws.Range("L15").Value = 0.3583950354711902
$ws.Range("N15").Value = 1.399984877553569
$ws.Range("O15").Value = 2.274345432317489
$ws.Range("B16").Value = 0.7145781811015581
$ws.Range("C16").Value = 0.2690359820902311
$ws.Range("D16").Value = 0.04344036383249517
$ws.Range("F16").Value = 0.6763587007867287
$ws.Range("G16").Value = 0.5097503013812172
$ws.Range("H16").Value = 0.6459115575200798
$ws.Range("K16").Value = 0.3973584120948601
$ws.Range("L16").Value = 0.3487718704368206
$ws.Range("N16").Value = 1.409278591853866
$ws.Range("O16").Value = 2.276725506235323
$ws.Range("B17").Value = 0.6933709610334802
$ws.Range("C17").Value = 0.269912494166733
$ws.Range("D17").Value = 0.04220497000907386
$ws.Range("F17").Value = 0.6751925090319233
$ws.Range("G17").Value = 0.5094815895110258
$ws.Range("H17").Value = 0.6470802287291662
$ws.Range("K17").Value = 0.381431632526386
$ws.Range("L17").Value = 0.3429140319069717
$ws.Range("N17").Value = 1.415108624245134
$ws.Range("O17").Value = 2.27851145019639
$ws.Range("B18").Value = 0.681184742605069
$ws.Range("C18").Value = 0.270423957476047
$ws.Range("D18").Value = 0.04149294651789859
$ws.Range("F18").Value = 0.6745595600624696
$ws.Range("G18").Value = 0.5093586605434126
$ws.Range("H18").Value = 0.6477779860893094
$ws.Range("K18").Value = 0.3722647429362098
$ws.Range("L18").Value = 0.3395615726655592
$ws.Range("N18").Value = 1.418509169694362
$ws.Range("O18").Value = 2.279658446940601
$ws.Range("B19").Value = 0.6770607173334326
$ws.Range("C19").Value = 0.2705983873991222
$ws.Range("D19").Value = 0.04125161841108138
$ws.Range("F19").Value = 0.6743517508219128
$ws.Range("G19").Value = 0.5093224702354036
$ws.Range("H19").Value = 0.6480186278752171
$ws.Range("K19").Value = 0.3691599398486574
$ws.Range("L19").Value = 0.3384293787003116
$ws.Range("N19").Value = 1.4196686528781
$ws.Range("O19").Value = 2.28006737282621
$ws.Range("B20").Value = 0.6956273096090797
$ws.Range("C20").Value = 0.2698184308411964
$ws.Range("D20").Value = 0.04233663092641393
$ws.Range("F20").Value = 0.6753127390056051
$ws.Range("G20").Value = 0.5095069209218934
$ws.Range("H20").Value = 0.6469531758352929
$ws.Range("K20").Value = 0.3831277141158012
$ws.Range("L20").Value = 0.3435358693427162
$ws.Range("N20").Value = 1.414483115378255
$ws.Range("O20").Value = 2.27830893825697
$ws.Range("B21").Value = 0.7581454796759317
$ws.Range("C21").Value = 0.2672852603461138
$ws.Range("D21").Value = 0.04596462255916833
$ws.Range("F21").Value = 0.678993179597029
$ws.Range("G21").Value = 0.5105045700779414
$ws.Range("H21").Value = 0.643678766101047
$ws.Range("K21").Value = 0.4299818072478558
$ws.Range("L21").Value = 0.3608936363240502
$ws.Range("N21").Value = 1.397625201577597
$ws.Range("O21").Value = 2.273831664652135
$ws.Range("B22").Value = 0.7991026997244148
$ws.Range("C22").Value = 0.2656947137683217
$ws.Range("D22").Value = 0.04832248338753686
$ws.Range("F22").Value = 0.6817348630246087
$ws.Range("G22").Value = 0.511438364981359
$ws.Range("H22").Value = 0.6417660549005717
$ws.Range("K22").Value = 0.4605442355366165
$ws.Range("L22").Value = 0.3723869444373804
$ws.Range("N22").Value = 1.387028679997968
$ws.Range("O22").Value = 2.271971569837376
$ws.Range("B23").Value = 0.777234358623474
$ws.Range("C23").Value = 0.2665376745984345
$ws.Range("D23").Value = 0.04706527296191609
$ws.Range("F23").Value = 0.6802408239683331
$ws.Range("G23").Value = 0.5109141922153242
$ws.Range("H23").Value = 0.6427661248352337
$ws.Range("K23").Value = 0.4442381228382146
$ws.Range("L23").Value = 0.3662391512198724
$ws.Range("N23").Value = 1.392645701729363
$ws.Range("O23").Value = 2.2728668155425
$ws.Range("B24").Value = 0.694607195373294
$ws.Range("C24").Value = 0.2698609333685118
$ws.Range("D24").Value = 0.04227711256464772
$ws.Range("F24").Value = 0.6752582661890827
$ws.Range("G24").Value = 0.5094953702948715
$ws.Range("H24").Value = 0.6470105358527292
$ws.Range("K24").Value = 0.38236094789238
$ws.Range("L24").Value = 0.3432546890046524
$ws.Range("N24").Value = 1.414765755985361
$ws.Range("O24").Value = 2.278400119380336
$ws.Range("B25").Value = 0.6059273098708786
$ws.Range("C25").Value = 0.2737254265501647
$ws.Range("D25").Value = 0.03705671340148342
$ws.Range("F25").Value = 0.6713296575144838
$ws.Range("G25").Value = 0.5091739849771102
$ws.Range("H25").Value = 0.6525665849687385
$ws.Range("K25").Value = 0.3153804308811061
$ws.Range("L25").Value = 0.3191053919608748
$ws.Range("N25").Value = 1.440432804138551
$ws.Range("O25").Value = 2.288947890832944

